# Began working on framework model and using new standards
#
# Adds two new rows to the Sheet1 parameter table (inside Table1):
#   StartPopUpText  / <email>                         / "The popup text to be displayed..."
#   StartPopUpTitle / "Companies House Demo Robot"     / "The title of the start pop up."
# inserted just above the existing "WorkpackageName" row, and expands the
# table / sheet dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The parameter table currently looks like (1-indexed sheet rows):
#   1  Name              Value                                 Description
#   2   ROBOT PARAMETERS
#   3  boolBreakpoint1    TRUE                                  breakpoint parameter
#   4  LogFilePath        C:\Users\{0}\...xlsx                  log file path
#   5  strEmailAccount    Sean.Crotty@defra.gov.uk               outlook email account
#   6  WorkpackageName    Demo Robot                            workpackage name
#   7   ROBOT PARAMETERS (2nd section label = "Regex")
#   8  (blank templated row)
#
# Insert two new blank rows above row 6 (the "WorkpackageName" row), which
# shifts it (and everything below) down to rows 8/9/10, then fill the two
# freed-up rows with the new StartPopUpText / StartPopUpTitle parameters.

$ws.Range("A6:C7").Insert()

$ws.Range("A6").Value = "StartPopUpText"
$ws.Range("B6").Value = "Sean.Crotty@defra.gov.uk"
$ws.Range("C6").Value = "The popup text to be displayed to the user at the start of the process."
$ws.Rows.Item(6).RowHeight = 30

$ws.Range("A7").Value = "StartPopUpTitle"
$ws.Range("B7").Value = "Companies House Demo Robot"
$ws.Range("C7").Value = "The title of the start pop up."

# Grow Table1 so it keeps covering the whole parameter block (was A1:C8).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C10"))

# Match the saved selection recorded in the edited workbook.
[void]$ws.Range("D7").Select()
